$d = $word.ActiveDocument

# --- "Process" list paragraph ---
# "sustain, discontinue, take off, Set up, Haste into, refine"
# ->  "sustain, take off, discontinue, set up, haste into, refine"
$pProcess = $d.Paragraphs(8).Range
$pProcess.Find.Execute("discontinue, take off, Set up", $true, $false, $false, $false, $false, $true, 1, $false, "take off, discontinue, set up", 2) | Out-Null

$pProcess2 = $d.Paragraphs(8).Range
$pProcess2.Find.Execute("Haste into", $true, $false, $false, $false, $false, $true, 1, $false, "haste into", 2) | Out-Null

# --- "Outcome" list paragraph ---
# "Arrive at, unforeseen, take for granted"
# -> "scrutinize, arrive at, unforeseen, take for granted"
$pOutcome = $d.Paragraphs(11).Range
$pOutcome.Find.Execute("Arrive at", $true, $false, $false, $false, $false, $true, 1, $false, "scrutinize, arrive at", 2) | Out-Null

# --- "Product" list paragraph ---
# "scale down, desirable, scrutinize, lose sight of, upmarket, withhold"
# -> "scale down, desirable, lose sight of, upmarket, withhold"
$pProduct = $d.Paragraphs(14).Range
$pProduct.Find.Execute("desirable, scrutinize", $true, $false, $false, $false, $false, $true, 1, $false, "desirable", 2) | Out-Null

# the "Product" heading run/paragraph mark picked up an explicit en-US language
# tag as a side effect of the original edit; reproduce that on the run.
$d.Paragraphs(13).Range.Font.LanguageID = "en-US"
